# Auto-generated edit script applying numeric corrections to the
# 'currentAveragePrice*' / 'LevePrice*' / 'LeveProfit*' columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1184.1765

$ws.Range("H32").Value = 5908.8184
$ws.Range("I32").Value = 2333.3333
$ws.Range("J32").Value = 10199.4
$ws.Range("K32").Value = 2333.3333
$ws.Range("L32").Value = 10199.4
$ws.Range("M32").Value = -2007.3333
$ws.Range("N32").Value = -10851.4

$ws.Range("H33").Value = 147.09091
$ws.Range("I33").Value = 146
$ws.Range("K33").Value = 146
$ws.Range("M33").Value = 83

$ws.Range("H74").Value = 5964.2144
$ws.Range("I74").Value = 5964.2144
$ws.Range("K74").Value = 5964.2144
$ws.Range("M74").Value = -5028.2144

$ws.Range("H77").Value = 5964.2144
$ws.Range("I77").Value = 5964.2144
$ws.Range("K77").Value = 29821.072
$ws.Range("M77").Value = -25141.072

$ws.Range("H111").Value = 300
$ws.Range("I111").Value = 300
$ws.Range("K111").Value = 900
$ws.Range("M111").Value = 2167

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H125").Value = 2438.5715
$ws.Range("I125").Value = 1759.8
$ws.Range("K125").Value = 15838.2
$ws.Range("M125").Value = -13378.2

$ws.Range("H137").Value = 6159
$ws.Range("J137").Value = 18017.166
$ws.Range("L137").Value = 54051.49800000001
$ws.Range("N137").Value = -59151.49800000001

$ws.Range("H138").Value = 2836.3606
$ws.Range("J138").Value = 2529.3333
$ws.Range("L138").Value = 7587.999899999999
$ws.Range("N138").Value = -17867.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 245319.95
$ws.Range("I32").Value = 274708.66
$ws.Range("J32").Value = 27843.6
$ws.Range("K32").Value = 274708.66
$ws.Range("L32").Value = 27843.6
$ws.Range("M32").Value = -274421.66
$ws.Range("N32").Value = -28417.6

$ws.Range("H45").Value = 3937.25
$ws.Range("J45").Value = 4108
$ws.Range("L45").Value = 4108
$ws.Range("N45").Value = -4862

$ws.Range("H61").Value = 3812.4285
$ws.Range("I61").Value = 3525.875
$ws.Range("K61").Value = 3525.875
$ws.Range("M61").Value = -3313.875

$ws.Range("H63").Value = 3860
$ws.Range("I63").Value = 3795
$ws.Range("J63").Value = 3990
$ws.Range("K63").Value = 3795
$ws.Range("L63").Value = 3990
$ws.Range("M63").Value = -3109
$ws.Range("N63").Value = -5362

$ws.Range("H66").Value = 3860
$ws.Range("I66").Value = 3795
$ws.Range("J66").Value = 3990
$ws.Range("K66").Value = 18975
$ws.Range("L66").Value = 19950
$ws.Range("M66").Value = -15543
$ws.Range("N66").Value = -26814

$ws.Range("H74").Value = 5130.0166
$ws.Range("I74").Value = 3347.0908
$ws.Range("J74").Value = 9744.647000000001
$ws.Range("K74").Value = 3347.0908
$ws.Range("L74").Value = 9744.647000000001
$ws.Range("M74").Value = -2473.0908
$ws.Range("N74").Value = -11492.647

$ws.Range("H77").Value = 5130.0166
$ws.Range("I77").Value = 3347.0908
$ws.Range("J77").Value = 9744.647000000001
$ws.Range("K77").Value = 16735.454
$ws.Range("L77").Value = 48723.235
$ws.Range("M77").Value = -12367.454
$ws.Range("N77").Value = -57459.235

$ws.Range("H97").Value = 1103.2142
$ws.Range("I97").Value = 937.0833
$ws.Range("J97").Value = 2100
$ws.Range("K97").Value = 937.0833
$ws.Range("L97").Value = 2100
$ws.Range("M97").Value = -441.0833
$ws.Range("N97").Value = -3092

$ws.Range("H132").Value = 5425.3228
$ws.Range("I132").Value = 3877.7026
$ws.Range("K132").Value = 11633.1078
$ws.Range("M132").Value = -9103.1078

$ws.Range("H136").Value = 3812.4285
$ws.Range("I136").Value = 3525.875
$ws.Range("K136").Value = 10577.625
$ws.Range("M136").Value = -8027.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 51193.477
$ws.Range("I20").Value = 80919
$ws.Range("K20").Value = 80919
$ws.Range("M20").Value = -80672

$ws.Range("H80").Value = 630.96155
$ws.Range("J80").Value = 593
$ws.Range("L80").Value = 593
$ws.Range("N80").Value = -2589

$ws.Range("H83").Value = 630.96155
$ws.Range("J83").Value = 593
$ws.Range("L83").Value = 2965
$ws.Range("N83").Value = -12949

$ws.Range("H86").Value = 1986.1578
$ws.Range("I86").Value = 1843
$ws.Range("J86").Value = 2749.6667
$ws.Range("K86").Value = 1843
$ws.Range("L86").Value = 2749.6667
$ws.Range("M86").Value = -720
$ws.Range("N86").Value = -4995.6667

$ws.Range("H89").Value = 1986.1578
$ws.Range("I89").Value = 1843
$ws.Range("J89").Value = 2749.6667
$ws.Range("K89").Value = 9215
$ws.Range("L89").Value = 13748.3335
$ws.Range("M89").Value = -3599
$ws.Range("N89").Value = -24980.3335

$ws.Range("H107").Value = 1716.8948
$ws.Range("I107").Value = 1683.5883
$ws.Range("K107").Value = 1683.5883
$ws.Range("M107").Value = 236.4117000000001

$ws.Range("H134").Value = 3242.0454
$ws.Range("I134").Value = 3640.3333
$ws.Range("J134").Value = 1449.75
$ws.Range("K134").Value = 10920.9999
$ws.Range("L134").Value = 4349.25
$ws.Range("M134").Value = -8385.999899999999
$ws.Range("N134").Value = -9419.25

$ws.Range("H137").Value = 56391.168
$ws.Range("J137").Value = 68099.14
$ws.Range("L137").Value = 68099.14
$ws.Range("N137").Value = -78299.14

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2321
$ws.Range("I22").Value = 1100.6
$ws.Range("J22").Value = 2999
$ws.Range("K22").Value = 1100.6
$ws.Range("L22").Value = 2999
$ws.Range("M22").Value = -750.5999999999999
$ws.Range("N22").Value = -3699

$ws.Range("H31").Value = 2840.75
$ws.Range("I31").Value = 4398
$ws.Range("K31").Value = 4398
$ws.Range("M31").Value = -4103

$ws.Range("H34").Value = 2840.75
$ws.Range("I34").Value = 4398
$ws.Range("K34").Value = 4398
$ws.Range("M34").Value = -4196

$ws.Range("H58").Value = 2673.4595
$ws.Range("I58").Value = 2530.0952
$ws.Range("J58").Value = 2861.625
$ws.Range("K58").Value = 2530.0952
$ws.Range("L58").Value = 2861.625
$ws.Range("M58").Value = -2327.0952
$ws.Range("N58").Value = -3267.625

$ws.Range("H62").Value = 2505.2
$ws.Range("J62").Value = 2218.6667
$ws.Range("L62").Value = 2218.6667
$ws.Range("N62").Value = -3466.6667

$ws.Range("H65").Value = 2505.2
$ws.Range("J65").Value = 2218.6667
$ws.Range("L65").Value = 11093.3335
$ws.Range("N65").Value = -17333.3335

$ws.Range("H86").Value = 98266.5
$ws.Range("J86").Value = 15000
$ws.Range("L86").Value = 15000
$ws.Range("N86").Value = -17246

$ws.Range("H89").Value = 98266.5
$ws.Range("J89").Value = 15000
$ws.Range("L89").Value = 75000
$ws.Range("N89").Value = -86232

$ws.Range("H99").Value = 13728
$ws.Range("I99").Value = 35622.832
$ws.Range("K99").Value = 35622.832
$ws.Range("M99").Value = -34124.832

$ws.Range("H105").Value = 13480.875
$ws.Range("I105").Value = 15321.143
$ws.Range("J105").Value = 599
$ws.Range("K105").Value = 15321.143
$ws.Range("L105").Value = 599
$ws.Range("M105").Value = -13574.143
$ws.Range("N105").Value = -4093

$ws.Range("H107").Value = 809.9231
$ws.Range("I107").Value = 590.35297
$ws.Range("K107").Value = 590.35297
$ws.Range("M107").Value = 1329.64703

$ws.Range("H126").Value = 13728
$ws.Range("I126").Value = 35622.832
$ws.Range("K126").Value = 106868.496
$ws.Range("M126").Value = -104398.496

$ws.Range("H132").Value = 2369.1226
$ws.Range("I132").Value = 2168.6667
$ws.Range("K132").Value = 6506.000100000001
$ws.Range("M132").Value = -3976.000100000001

$ws.Range("H136").Value = 2673.4595
$ws.Range("I136").Value = 2530.0952
$ws.Range("J136").Value = 2861.625
$ws.Range("K136").Value = 7590.285600000001
$ws.Range("L136").Value = 8584.875
$ws.Range("M136").Value = -5040.285600000001
$ws.Range("N136").Value = -13684.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1075.258
$ws.Range("I5").Value = 746.05
$ws.Range("J5").Value = 1673.8182
$ws.Range("K5").Value = 2238.15
$ws.Range("L5").Value = 5021.4546
$ws.Range("M5").Value = -2126.15
$ws.Range("N5").Value = -5245.4546

$ws.Range("H6").Value = 1700
$ws.Range("I6").Value = 750
$ws.Range("K6").Value = 2250
$ws.Range("M6").Value = -2137

$ws.Range("H11").Value = 166667070
$ws.Range("I11").Value = 550
$ws.Range("K11").Value = 1650
$ws.Range("M11").Value = -1510

$ws.Range("H97").Value = 674.75
$ws.Range("I97").Value = 749.5
$ws.Range("J97").Value = 600
$ws.Range("K97").Value = 2248.5
$ws.Range("L97").Value = 1800
$ws.Range("M97").Value = -1752.5
$ws.Range("N97").Value = -2792

$ws.Range("H107").Value = 2913
$ws.Range("I107").Value = 788.625
$ws.Range("J107").Value = 4458
$ws.Range("K107").Value = 2365.875
$ws.Range("L107").Value = 13374
$ws.Range("M107").Value = -445.875
$ws.Range("N107").Value = -17214

$ws.Range("H122").Value = 949721.5
$ws.Range("J122").Value = 1271.5834
$ws.Range("L122").Value = 11444.2506
$ws.Range("N122").Value = -16344.2506

$ws.Range("H123").Value = 11499.833
$ws.Range("J123").Value = 15000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -49900

$ws.Range("H135").Value = 1075.258
$ws.Range("I135").Value = 746.05
$ws.Range("J135").Value = 1673.8182
$ws.Range("K135").Value = 6714.45
$ws.Range("L135").Value = 15064.3638
$ws.Range("M135").Value = -4179.45
$ws.Range("N135").Value = -20134.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 9999.5
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H14").Value = 75003080
$ws.Range("I14").Value = 81820820
$ws.Range("J14").Value = 8000
$ws.Range("K14").Value = 81820820
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = -81820652
$ws.Range("N14").Value = -8336

$ws.Range("H53").Value = 19993
$ws.Range("I53").Value = 19999
$ws.Range("J53").Value = 19990
$ws.Range("K53").Value = 19999
$ws.Range("L53").Value = 19990
$ws.Range("N53").Value = -21252
$ws.Range("M53").Value = -19368

$ws.Range("H80").Value = 2748.125
$ws.Range("I80").Value = 2672.25
$ws.Range("K80").Value = 2672.25
$ws.Range("M80").Value = -1674.25

$ws.Range("H83").Value = 2748.125
$ws.Range("I83").Value = 2672.25
$ws.Range("K83").Value = 13361.25
$ws.Range("M83").Value = -8369.25

$ws.Range("H102").Value = 3056.238
$ws.Range("I102").Value = 3346.6875
$ws.Range("J102").Value = 2126.8
$ws.Range("K102").Value = 3346.6875
$ws.Range("L102").Value = 2126.8
$ws.Range("M102").Value = -1724.6875
$ws.Range("N102").Value = -5370.8

$ws.Range("H107").Value = 529.3333
$ws.Range("J107").Value = 550.7143
$ws.Range("L107").Value = 550.7143
$ws.Range("N107").Value = -4390.7143

$ws.Range("H113").Value = 1712.1666
$ws.Range("I113").Value = 1616.2222
$ws.Range("K113").Value = 1616.2222
$ws.Range("M113").Value = 553.7778000000001

$ws.Range("H122").Value = 3232.05
$ws.Range("I122").Value = 3110.1875
$ws.Range("J122").Value = 3719.5
$ws.Range("K122").Value = 9330.5625
$ws.Range("L122").Value = 11158.5
$ws.Range("M122").Value = -6880.5625
$ws.Range("N122").Value = -16058.5

$ws.Range("H126").Value = 2621.6316
$ws.Range("I126").Value = 2466
$ws.Range("K126").Value = 7398
$ws.Range("M126").Value = -4928

$ws.Range("H132").Value = 8549.3125
$ws.Range("I132").Value = 12253.223
$ws.Range("K132").Value = 36759.669
$ws.Range("M132").Value = -34229.669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9000
$ws.Range("I7").Value = 8000
$ws.Range("K7").Value = 8000
$ws.Range("M7").Value = -7888

$ws.Range("H34").Value = 11000
$ws.Range("I34").Value = 11000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 11000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -10828
$ws.Range("N34").ClearContents()

$ws.Range("H40").Value = 4716.643
$ws.Range("I40").Value = 4171.5
$ws.Range("J40").Value = 7987.5
$ws.Range("K40").Value = 4171.5
$ws.Range("L40").Value = 7987.5
$ws.Range("M40").Value = -4035.5
$ws.Range("N40").Value = -8259.5

$ws.Range("H55").Value = 1495.6875
$ws.Range("I55").Value = 1414.8462
$ws.Range("J55").Value = 1551
$ws.Range("K55").Value = 1414.8462
$ws.Range("L55").Value = 1551
$ws.Range("M55").Value = -1241.8462
$ws.Range("N55").Value = -1897

$ws.Range("H61").Value = 12944.85
$ws.Range("I61").Value = 12009.467
$ws.Range("K61").Value = 12009.467
$ws.Range("M61").Value = -11807.467

$ws.Range("H68").Value = 7078
$ws.Range("J68").Value = 4931.3335
$ws.Range("L68").Value = 4931.3335
$ws.Range("N68").Value = -6429.3335

$ws.Range("H71").Value = 7078
$ws.Range("J71").Value = 4931.3335
$ws.Range("L71").Value = 24656.6675
$ws.Range("N71").Value = -32144.6675

$ws.Range("H100").Value = 2133.7
$ws.Range("I100").Value = 2133.7
$ws.Range("K100").Value = 2133.7
$ws.Range("M100").Value = -1592.7

$ws.Range("H113").Value = 12944.85
$ws.Range("I113").Value = 12009.467
$ws.Range("K113").Value = 12009.467
$ws.Range("M113").Value = -9839.467000000001

$ws.Range("H122").Value = 4762.3335
$ws.Range("I122").Value = 3717.8
$ws.Range("K122").Value = 11153.4
$ws.Range("M122").Value = -8703.400000000001

$ws.Range("H126").Value = 9000
$ws.Range("I126").Value = 8000
$ws.Range("K126").Value = 24000
$ws.Range("M126").Value = -21530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8495.883
$ws.Range("I62").Value = 9168.429
$ws.Range("J62").Value = 8025.1
$ws.Range("K62").Value = 9168.429
$ws.Range("L62").Value = 8025.1
$ws.Range("M62").Value = -8544.429
$ws.Range("N62").Value = -9273.1

$ws.Range("H65").Value = 8495.883
$ws.Range("I65").Value = 9168.429
$ws.Range("J65").Value = 8025.1
$ws.Range("K65").Value = 45842.145
$ws.Range("L65").Value = 40125.5
$ws.Range("M65").Value = -42722.145
$ws.Range("N65").Value = -46365.5

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H96").Value = 1151.9166
$ws.Range("I96").Value = 1398
$ws.Range("K96").Value = 1398
$ws.Range("M96").Value = -25

$ws.Range("H107").Value = 1947.0667
$ws.Range("I107").Value = 1104.579
$ws.Range("J107").Value = 3402.2727
$ws.Range("K107").Value = 3313.737
$ws.Range("L107").Value = 10206.8181
$ws.Range("M107").Value = -1393.737
$ws.Range("N107").Value = -14046.8181

$ws.Range("H117").Value = 25204.5
$ws.Range("J117").Value = 25204.5
$ws.Range("L117").Value = 25204.5
$ws.Range("N117").Value = -34382.5

$ws.Range("H122").Value = 91200.62
$ws.Range("I122").Value = 4835.875
$ws.Range("K122").Value = 14507.625
$ws.Range("M122").Value = -12057.625

$ws.Range("H126").Value = 2591.8462
$ws.Range("I126").Value = 2369.5
$ws.Range("K126").Value = 7108.5
$ws.Range("M126").Value = -4638.5

$ws.Range("H132").Value = 4072.75

$ws.Range("H136").Value = 1507.8695
$ws.Range("I136").Value = 1180.25
$ws.Range("K136").Value = 3540.75
$ws.Range("M136").Value = -990.75

